$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (rows 2-97) from 45183 to 45184 (date serial +1 day)
for ($row = 2; $row -le 97; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45183) {
        $cell.Value2 = 45184
    }
}
